$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 with flight data (mirrors row 10's flight reusing date "Sunday, Jan 15")
$ws.Range("A11").Value = 10.0
$ws.Range("B11").Value = "Sunday, Jan 15"
$ws.Range("C11").Value = "4:35 PM"
$ws.Range("D11").Value = "FR2679"
$ws.Range("E11").Value = "London"
$ws.Range("F11").Value = "(STN)"
$ws.Range("G11").Value = "Ryanair "
$ws.Range("H11").Value = "B738"
$ws.Range("I11").Value = "(EI-DYN)"
$ws.Range("J11").Value = "5:18 PM"
$ws.Range("L11").Value = "0 hours, 43 minutes"

# Materialize the blank K11/M11 cells (present in all other data rows) without
# introducing a new cell style, matching the existing blank cells elsewhere.
$ws.Range("K11").Borders.LineStyle = -4142
$ws.Range("M11").Borders.LineStyle = -4142
